$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
[void]$ws.Select()

# Add "ok" marker cells in column C for the steps that should show in the
# architect view (rows 34, 35, 36 and 38 - matching existing rows 31-33).
$ws.Range("C34").Value = "ok"
$ws.Range("C35").Value = "ok"
$ws.Range("C36").Value = "ok"
$ws.Range("C38").Value = "ok"

# Add the new "Aparecer en vista arq" note next to the steps that already
# involve the architect profile (rows 37 and 39).
$ws.Range("E37").Value = "Aparecer en vista arq"
$ws.Range("E39").Value = "Aparecer en vista arq"

# Update the view state to reflect where the user ended up scrolled to /
# selecting after making the edits above.
[void]$ws.Range("D48").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 27
$win.ScrollColumn = 1
